$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the tooltip text for "Fiscal Year End" (row 3, column E) to end with a period.
$ws.Range("E3").Value = "The date the fiscal year ends."

# Move/confirm the active selection to E3, matching the saved cursor position in the workbook.
$null = $ws.Range("E3").Select()
